$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 493.2857
$ws.Range("I19").Value = 309.94116
$ws.Range("J19").Value = 776.63635
$ws.Range("K19").Value = 309.94116
$ws.Range("L19").Value = 776.63635
$ws.Range("M19").Value = -134.94116
$ws.Range("N19").Value = -1126.63635
$ws.Range("H33").Value = 3872.1724
$ws.Range("I33").Value = 81.125
$ws.Range("J33").Value = 5316.381
$ws.Range("K33").Value = 81.125
$ws.Range("L33").Value = 5316.381
$ws.Range("M33").Value = 147.875
$ws.Range("N33").Value = -5774.381
$ws.Range("H51").Value = 4918.8
$ws.Range("J51").Value = 4918.8
$ws.Range("L51").Value = 4918.8
$ws.Range("N51").Value = -5886.8
$ws.Range("H74").Value = 7043.4287
$ws.Range("I74").Value = 8042.857
$ws.Range("K74").Value = 8042.857
$ws.Range("M74").Value = -7106.857
$ws.Range("H75").Value = 63000
$ws.Range("J75").Value = 63000
$ws.Range("L75").Value = 63000
$ws.Range("N75").Value = -64872
$ws.Range("H77").Value = 7043.4287
$ws.Range("I77").Value = 8042.857
$ws.Range("K77").Value = 40214.285
$ws.Range("M77").Value = -35534.285
$ws.Range("H78").Value = 63000
$ws.Range("J78").Value = 63000
$ws.Range("L78").Value = 189000
$ws.Range("N78").Value = -198360
$ws.Range("H98").Value = 4675434.5
$ws.Range("I98").Value = 7010508.5
$ws.Range("J98").Value = 5286.25
$ws.Range("K98").Value = 7010508.5
$ws.Range("L98").Value = 5286.25
$ws.Range("M98").Value = -7009010.5
$ws.Range("N98").Value = -8282.25
$ws.Range("H100").Value = 2119.2856
$ws.Range("I100").Value = 1500.4166
$ws.Range("J100").Value = 2944.4443
$ws.Range("K100").Value = 1500.4166
$ws.Range("L100").Value = 2944.4443
$ws.Range("M100").Value = -959.4166
$ws.Range("N100").Value = -4026.4443
$ws.Range("H122").Value = 4675434.5
$ws.Range("I122").Value = 7010508.5
$ws.Range("J122").Value = 5286.25
$ws.Range("K122").Value = 21031525.5
$ws.Range("L122").Value = 15858.75
$ws.Range("M122").Value = -21029075.5
$ws.Range("N122").Value = -20758.75
$ws.Range("H137").Value = 906.34283
$ws.Range("I137").Value = 801.5454999999999
$ws.Range("J137").Value = 1083.6923
$ws.Range("K137").Value = 2404.6365
$ws.Range("L137").Value = 3251.0769
$ws.Range("M137").Value = 145.3635000000004
$ws.Range("N137").Value = -8351.0769
$ws.Range("H138").Value = 3095.4
$ws.Range("I138").Value = 1684.7667
$ws.Range("J138").Value = 3941.78
$ws.Range("K138").Value = 5054.300099999999
$ws.Range("L138").Value = 11825.34
$ws.Range("M138").Value = 85.69990000000053
$ws.Range("N138").Value = -22105.34

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 454039.12
$ws.Range("I32").Value = 4075.4363
$ws.Range("J32").Value = 3203817.2
$ws.Range("K32").Value = 4075.4363
$ws.Range("L32").Value = 3203817.2
$ws.Range("M32").Value = -3788.4363
$ws.Range("N32").Value = -3204391.2
$ws.Range("H61").Value = 1412.6666
$ws.Range("I61").Value = 1403
$ws.Range("J61").Value = 1432
$ws.Range("K61").Value = 1403
$ws.Range("L61").Value = 1432
$ws.Range("M61").Value = -1191
$ws.Range("N61").Value = -1856
$ws.Range("H74").Value = 1014.2857
$ws.Range("I74").Value = 964.3200000000001
$ws.Range("J74").Value = 1139.2
$ws.Range("K74").Value = 964.3200000000001
$ws.Range("L74").Value = 1139.2
$ws.Range("M74").Value = -90.32000000000005
$ws.Range("N74").Value = -2887.2
$ws.Range("H77").Value = 1014.2857
$ws.Range("I77").Value = 964.3200000000001
$ws.Range("J77").Value = 1139.2
$ws.Range("K77").Value = 4821.6
$ws.Range("L77").Value = 5696
$ws.Range("M77").Value = -453.6000000000004
$ws.Range("N77").Value = -14432
$ws.Range("H132").Value = 21299532
$ws.Range("I132").Value = 25642438
$ws.Range("K132").Value = 76927314
$ws.Range("M132").Value = -76924784
$ws.Range("H136").Value = 1412.6666
$ws.Range("I136").Value = 1403
$ws.Range("J136").Value = 1432
$ws.Range("K136").Value = 4209
$ws.Range("L136").Value = 4296
$ws.Range("M136").Value = -1659
$ws.Range("N136").Value = -9396

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1610.2858
$ws.Range("I99").Value = 987.06665
$ws.Range("J99").Value = 3168.3333
$ws.Range("K99").Value = 987.06665
$ws.Range("L99").Value = 3168.3333
$ws.Range("M99").Value = 510.93335
$ws.Range("N99").Value = -6164.3333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1327.6111
$ws.Range("I16").Value = 1200.2
$ws.Range("J16").Value = 1486.875
$ws.Range("K16").Value = 1200.2
$ws.Range("L16").Value = 1486.875
$ws.Range("M16").Value = -913.2
$ws.Range("N16").Value = -2060.875
$ws.Range("H31").Value = 11907578
$ws.Range("I31").Value = 16669208
$ws.Range("J31").Value = 3502.3333
$ws.Range("K31").Value = 16669208
$ws.Range("L31").Value = 3502.3333
$ws.Range("M31").Value = -16668913
$ws.Range("N31").Value = -4092.3333
$ws.Range("H34").Value = 11907578
$ws.Range("I34").Value = 16669208
$ws.Range("J34").Value = 3502.3333
$ws.Range("K34").Value = 16669208
$ws.Range("L34").Value = 3502.3333
$ws.Range("M34").Value = -16669006
$ws.Range("N34").Value = -3906.3333
$ws.Range("H99").Value = 251429.8
$ws.Range("I99").Value = 370897.3
$ws.Range("J99").Value = 1634.091
$ws.Range("K99").Value = 370897.3
$ws.Range("L99").Value = 1634.091
$ws.Range("M99").Value = -369399.3
$ws.Range("N99").Value = -4630.091
$ws.Range("H113").Value = 1327.6111
$ws.Range("I113").Value = 1200.2
$ws.Range("J113").Value = 1486.875
$ws.Range("K113").Value = 1200.2
$ws.Range("L113").Value = 1486.875
$ws.Range("M113").Value = 969.8
$ws.Range("N113").Value = -5826.875
$ws.Range("H126").Value = 251429.8
$ws.Range("I126").Value = 370897.3
$ws.Range("J126").Value = 1634.091
$ws.Range("K126").Value = 1112691.9
$ws.Range("L126").Value = 4902.272999999999
$ws.Range("M126").Value = -1110221.9
$ws.Range("N126").Value = -9842.272999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 55557220
$ws.Range("I22").Value = 1333.3334
$ws.Range("J22").Value = 83335160
$ws.Range("K22").Value = 4000.0002
$ws.Range("L22").Value = 250005480
$ws.Range("M22").Value = -3831.0002
$ws.Range("N22").Value = -250005818
$ws.Range("H27").Value = 55557220
$ws.Range("I27").Value = 1333.3334
$ws.Range("J27").Value = 83335160
$ws.Range("K27").Value = 4000.0002
$ws.Range("L27").Value = 250005480
$ws.Range("M27").Value = -3898.0002
$ws.Range("N27").Value = -250005684
$ws.Range("H58").Value = 1393.3334
$ws.Range("J58").Value = 1590
$ws.Range("L58").Value = 4770
$ws.Range("N58").Value = -5026
$ws.Range("H97").Value = 358.1
$ws.Range("I97").Value = 619.75
$ws.Range("K97").Value = 1859.25
$ws.Range("M97").Value = -1363.25
$ws.Range("H137").Value = 1658.4117
$ws.Range("I137").Value = 1679.0769
$ws.Range("J137").Value = 1591.25
$ws.Range("K137").Value = 5037.2307
$ws.Range("L137").Value = 4773.75
$ws.Range("M137").Value = 62.76929999999993
$ws.Range("N137").Value = -14973.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 43208.2
$ws.Range("I58").Value = 67347
$ws.Range("K58").Value = 67347
$ws.Range("M58").Value = -67070
$ws.Range("H102").Value = 2430.923
$ws.Range("I102").Value = 1798.3334
$ws.Range("J102").Value = 2973.1428
$ws.Range("K102").Value = 1798.3334
$ws.Range("L102").Value = 2973.1428
$ws.Range("M102").Value = -176.3334
$ws.Range("N102").Value = -6217.1428
$ws.Range("H126").Value = 11913126
$ws.Range("I126").Value = 12142.444
$ws.Range("J126").Value = 33334896
$ws.Range("K126").Value = 36427.33199999999
$ws.Range("L126").Value = 100004688
$ws.Range("M126").Value = -33957.33199999999
$ws.Range("N126").Value = -100009628

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1935.0938
$ws.Range("I40").Value = 1719.9
$ws.Range("J40").Value = 2293.75
$ws.Range("K40").Value = 1719.9
$ws.Range("L40").Value = 2293.75
$ws.Range("M40").Value = -1583.9
$ws.Range("N40").Value = -2565.75
$ws.Range("H55").Value = 561.7037
$ws.Range("I55").Value = 291.93332
$ws.Range("J55").Value = 898.9167
$ws.Range("K55").Value = 291.93332
$ws.Range("L55").Value = 898.9167
$ws.Range("M55").Value = -118.93332
$ws.Range("N55").Value = -1244.9167
$ws.Range("H57").Value = 313995
$ws.Range("I57").Value = 600000
$ws.Range("J57").Value = 27990
$ws.Range("K57").Value = 600000
$ws.Range("L57").Value = 27990
$ws.Range("M57").Value = -599434
$ws.Range("N57").Value = -29122
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H68").Value = 1882.5122
$ws.Range("I68").Value = 1780
$ws.Range("J68").Value = 2060.2
$ws.Range("K68").Value = 1780
$ws.Range("L68").Value = 2060.2
$ws.Range("M68").Value = -1031
$ws.Range("N68").Value = -3558.2
$ws.Range("H71").Value = 1882.5122
$ws.Range("I71").Value = 1780
$ws.Range("J71").Value = 2060.2
$ws.Range("K71").Value = 8900
$ws.Range("L71").Value = 10301
$ws.Range("M71").Value = -5156
$ws.Range("N71").Value = -17789
$ws.Range("H132").Value = 287526.84
$ws.Range("I132").Value = 334788.47
$ws.Range("J132").Value = 3957.2
$ws.Range("K132").Value = 1004365.41
$ws.Range("L132").Value = 11871.6
$ws.Range("M132").Value = -1001835.41
$ws.Range("N132").Value = -16931.6
$ws.Range("H136").Value = 6404.0645
$ws.Range("I136").Value = 6474.409
$ws.Range("J136").Value = 6232.1113
$ws.Range("K136").Value = 19423.227
$ws.Range("L136").Value = 18696.3339
$ws.Range("M136").Value = -16873.227
$ws.Range("N136").Value = -23796.3339

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 22495
$ws.Range("J92").Value = 22495
$ws.Range("L92").Value = 22495
$ws.Range("N92").Value = -27487
